$d = $word.ActiveDocument

$replacements = @(
    @("803÷8=100, 3", "915÷5=183, 0"),
    @("106÷6=17, 4", "110÷3=36, 2"),
    @("110÷4=27, 2", "486÷7=69, 3"),
    @("804÷3=268, 0", "482÷8=60, 2"),
    @("313÷7=44, 5", "681÷5=136, 1"),
    @("708÷6=118, 0", "104÷8=13, 0"),
    @("795÷8=99, 3", "599÷8=74, 7"),
    @("474÷8=59, 2", "913÷4=228, 1"),
    @("226÷5=45, 1", "482÷5=96, 2"),
    @("511÷7=73, 0", "400÷3=133, 1"),
    @("163÷4=40, 3", "659÷3=219, 2"),
    @("183÷2=91, 1", "468÷4=117, 0"),
    @("691÷2=345, 1", "252÷9=28, 0"),
    @("267÷9=29, 6", "120÷7=17, 1"),
    @("579÷6=96, 3", "778÷7=111, 1"),
    @("744÷9=82, 6", "140÷5=28, 0"),
    @("792÷2=396, 0", "678÷2=339, 0"),
    @("569÷2=284, 1", "920÷2=460, 0"),
    @("167÷8=20, 7", "216÷9=24, 0"),
    @("123÷9=13, 6", "480÷3=160, 0"),
    @("680÷5=136, 0", "119÷4=29, 3"),
    @("526÷3=175, 1", "992÷6=165, 2"),
    @("712÷6=118, 4", "447÷5=89, 2"),
    @("332÷4=83, 0", "529÷9=58, 7"),
    @("888÷4=222, 0", "391÷8=48, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
